# Update "想去人数" (F column) values for the sheets that contain the
# full event listing data: "展览" and "全部类型" (sheet1 and sheet4).
# "演出" and "本地生活" only contain header rows, so they are skipped.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 11817
    3  = 11548
    4  = 612
    6  = 1042
    11 = 10859
    12 = 4191
    16 = 2475
    19 = 5
    21 = 457
    22 = 11161
    23 = 10958
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
